$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 900
$ws.Range("I18").Value = 900
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -616
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 234.8
$ws.Range("I33").Value = 188.95238
$ws.Range("J33").Value = 475.5
$ws.Range("K33").Value = 188.95238
$ws.Range("L33").Value = 475.5
$ws.Range("M33").Value = 40.04761999999999
$ws.Range("N33").Value = -933.5
$ws.Range("H40").Value = 3474215.5
$ws.Range("I40").Value = 6946033
$ws.Range("K40").Value = 6946033
$ws.Range("M40").Value = -6945858
$ws.Range("H74").Value = 7638
$ws.Range("I74").Value = 3846
$ws.Range("J74").Value = 8902
$ws.Range("K74").Value = 3846
$ws.Range("L74").Value = 8902
$ws.Range("M74").Value = -2910
$ws.Range("N74").Value = -10774
$ws.Range("H77").Value = 7638
$ws.Range("I77").Value = 3846
$ws.Range("J77").Value = 8902
$ws.Range("K77").Value = 19230
$ws.Range("L77").Value = 44510
$ws.Range("M77").Value = -14550
$ws.Range("N77").Value = -53870
$ws.Range("H100").Value = 6769.4
$ws.Range("I100").Value = 8914.615
$ws.Range("J100").Value = 2785.4285
$ws.Range("K100").Value = 8914.615
$ws.Range("L100").Value = 2785.4285
$ws.Range("M100").Value = -8373.615
$ws.Range("N100").Value = -3867.4285
$ws.Range("H116").Value = 2072.1667
$ws.Range("I116").Value = 1843.3334
$ws.Range("J116").Value = 2301
$ws.Range("K116").Value = 1843.3334
$ws.Range("L116").Value = 2301
$ws.Range("M116").Value = 1598.6666
$ws.Range("N116").Value = -9185
$ws.Range("H137").Value = 1440.6333
$ws.Range("I137").Value = 946.38464
$ws.Range("J137").Value = 1818.5883
$ws.Range("K137").Value = 2839.15392
$ws.Range("L137").Value = 5455.7649
$ws.Range("M137").Value = -289.1539199999997
$ws.Range("N137").Value = -10555.7649
$ws.Range("H138").Value = 2346.69
$ws.Range("J138").Value = 2964.0151
$ws.Range("L138").Value = 8892.0453
$ws.Range("N138").Value = -19172.0453

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6252.39
$ws.Range("I32").Value = 4980.7065
$ws.Range("J32").Value = 20876.75
$ws.Range("K32").Value = 4980.7065
$ws.Range("L32").Value = 20876.75
$ws.Range("M32").Value = -4693.7065
$ws.Range("N32").Value = -21450.75
$ws.Range("H74").Value = 10207154
$ws.Range("I74").Value = 12196148
$ws.Range("J74").Value = 13564.5
$ws.Range("K74").Value = 12196148
$ws.Range("L74").Value = 13564.5
$ws.Range("M74").Value = -12195274
$ws.Range("N74").Value = -15312.5
$ws.Range("H77").Value = 10207154
$ws.Range("I77").Value = 12196148
$ws.Range("J77").Value = 13564.5
$ws.Range("K77").Value = 60980740
$ws.Range("L77").Value = 67822.5
$ws.Range("M77").Value = -60976372
$ws.Range("N77").Value = -76558.5
$ws.Range("H102").Value = 1014.4286
$ws.Range("I102").Value = 1008.5
$ws.Range("J102").Value = 1050
$ws.Range("K102").Value = 1008.5
$ws.Range("L102").Value = 1050
$ws.Range("M102").Value = 613.5
$ws.Range("N102").Value = -4294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4874.5
$ws.Range("I20").Value = 4500
$ws.Range("J20").Value = 4999.3335
$ws.Range("K20").Value = 4500
$ws.Range("L20").Value = 4999.3335
$ws.Range("M20").Value = -4253
$ws.Range("N20").Value = -5493.3335
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H80").Value = 3920.7058
$ws.Range("I80").Value = 885.2308
$ws.Range("J80").Value = 5799.8096
$ws.Range("K80").Value = 885.2308
$ws.Range("L80").Value = 5799.8096
$ws.Range("M80").Value = 112.7692
$ws.Range("N80").Value = -7795.8096
$ws.Range("H83").Value = 3920.7058
$ws.Range("I83").Value = 885.2308
$ws.Range("J83").Value = 5799.8096
$ws.Range("K83").Value = 4426.154
$ws.Range("L83").Value = 28999.048
$ws.Range("M83").Value = 565.8459999999995
$ws.Range("N83").Value = -38983.048
$ws.Range("H99").Value = 2148.8125
$ws.Range("I99").Value = 1095
$ws.Range("J99").Value = 2299.3572
$ws.Range("K99").Value = 1095
$ws.Range("L99").Value = 2299.3572
$ws.Range("M99").Value = 403
$ws.Range("N99").Value = -5295.3572
$ws.Range("H132").Value = 59926.668
$ws.Range("J132").Value = 59926.668
$ws.Range("L132").Value = 59926.668
$ws.Range("N132").Value = -70046.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 450
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 100
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = 13
$ws.Range("N3").Value = -1026
$ws.Range("H22").Value = 20224.2
$ws.Range("I22").Value = 11360.223
$ws.Range("J22").Value = 100000
$ws.Range("K22").Value = 11360.223
$ws.Range("L22").Value = 100000
$ws.Range("M22").Value = -11010.223
$ws.Range("N22").Value = -100700
$ws.Range("H86").Value = 22775058
$ws.Range("I86").Value = 50051840
$ws.Range("K86").Value = 50051840
$ws.Range("M86").Value = -50050717
$ws.Range("H89").Value = 22775058
$ws.Range("I89").Value = 50051840
$ws.Range("K89").Value = 250259200
$ws.Range("M89").Value = -250253584
$ws.Range("H105").Value = 24372
$ws.Range("I105").Value = 36993.332
$ws.Range("J105").Value = 5440
$ws.Range("K105").Value = 36993.332
$ws.Range("L105").Value = 5440
$ws.Range("M105").Value = -35246.332
$ws.Range("N105").Value = -8934
$ws.Range("H107").Value = 489.14816
$ws.Range("I107").Value = 282.14285
$ws.Range("K107").Value = 282.14285
$ws.Range("M107").Value = 1637.85715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 124486600
$ws.Range("J32").Value = 124486600
$ws.Range("L32").Value = 373459800
$ws.Range("N32").Value = -373460366
$ws.Range("H113").Value = 10784935
$ws.Range("J113").Value = 7143504.5
$ws.Range("L113").Value = 21430513.5
$ws.Range("N113").Value = -21434853.5
$ws.Range("H131").Value = 816.7766
$ws.Range("I131").Value = 633.3333
$ws.Range("J131").Value = 836.2
$ws.Range("K131").Value = 1899.9999
$ws.Range("L131").Value = 2508.6
$ws.Range("M131").Value = 3140.0001
$ws.Range("N131").Value = -12588.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12502101
$ws.Range("I80").Value = 2474.75
$ws.Range("J80").Value = 25001726
$ws.Range("K80").Value = 2474.75
$ws.Range("L80").Value = 25001726
$ws.Range("M80").Value = -1476.75
$ws.Range("N80").Value = -25003722
$ws.Range("H83").Value = 12502101
$ws.Range("I83").Value = 2474.75
$ws.Range("J83").Value = 25001726
$ws.Range("K83").Value = 12373.75
$ws.Range("L83").Value = 125008630
$ws.Range("M83").Value = -7381.75
$ws.Range("N83").Value = -125018614
$ws.Range("H102").Value = 1228.5555
$ws.Range("I102").Value = 1183.1765
$ws.Range("K102").Value = 1183.1765
$ws.Range("M102").Value = 438.8235
$ws.Range("H113").Value = 2193.5625
$ws.Range("I113").Value = 1855.5555
$ws.Range("J113").Value = 2628.1428
$ws.Range("K113").Value = 1855.5555
$ws.Range("L113").Value = 2628.1428
$ws.Range("M113").Value = 314.4445000000001
$ws.Range("N113").Value = -6968.1428
$ws.Range("H126").Value = 4925
$ws.Range("I126").Value = 5850
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 17550
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -15080
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2956.0286
$ws.Range("I136").Value = 2734.3215
$ws.Range("J136").Value = 3842.8572
$ws.Range("K136").Value = 8202.9645
$ws.Range("L136").Value = 11528.5716
$ws.Range("M136").Value = -5652.9645
$ws.Range("N136").Value = -16628.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 62500708
$ws.Range("I81").Value = 62500708
$ws.Range("K81").Value = 125001416
$ws.Range("M81").Value = -125000355
$ws.Range("H84").Value = 62500708
$ws.Range("I84").Value = 62500708
$ws.Range("K84").Value = 625007080
$ws.Range("M84").Value = -625001776
$ws.Range("H113").Value = 47619450
$ws.Range("I113").Value = 55555820
$ws.Range("J113").Value = 1233.3334
$ws.Range("K113").Value = 166667460
$ws.Range("L113").Value = 3700.0002
$ws.Range("M113").Value = -166665290
$ws.Range("N113").Value = -8040.0002
